# "Move DevOps to separate file"
#
# The DevOps-related VPC/AZ diagram on slide 38 gets a few labels filled
# in with more specific (example) values:
#   - "Availability Zone"  -> "Availability Zone a"
#   - "region"             -> "region us-west-2"
#   - "VPC subnet"         -> "VPC subnet A"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(38)

# --- "Availability Zone" textbox (shape "TextBox 32") ---------------------
$azShape = $s.Shapes.Item(6)
$azRange = $azShape.TextFrame.TextRange
# Split the existing run right before "Zone" and type " a" after it, so the
# text becomes "Availability " + "Zone a" (mirrors editing the label in
# place in the PowerPoint UI).
$azZone = $azRange.Characters(14, 4)
$azZone.InsertAfter(" a") | Out-Null

# --- "region" textbox (shape "TextBox 33") --------------------------------
$regionShape = $s.Shapes.Item(7)
$regionShape.TextFrame.TextRange.Text = "region us-west-2"

# --- "VPC subnet" textbox (shape "TextBox 37") ----------------------------
$subnetShape = $s.Shapes.Item(9)
$subnetRange = $subnetShape.TextFrame.TextRange
$subnetWord = $subnetRange.Characters(5, 6)
$subnetWord.Text = "subnet A"
